# Refresh the crypto price/volume snapshot (GitHub Actions daily update).
# Rows 19/20 also swap rank order: ShibaInu moves to rank 19 (was Chainlink),
# Chainlink drops to rank 20 (was ShibaInu) - B/C/D/E updated accordingly.
#
# Note: several "Price" column values are plain digits-with-one-dot
# (e.g. 213.41) which Excel would otherwise auto-coerce to a number;
# a leading single-quote forces them to stay text, matching the
# original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.798.39'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '1.648.48'
$ws.Range("E3").Value = '  -0.46%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''213.41'
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''23.14'
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("D9").Value = '''0.259'
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").Value = '1.882.96'
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '1.641.92'
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("D15").Value = '''0.564'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '''64.45'
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("D17").Value = '27.778.68'
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").Value = '''233.60'
$ws.Range("E18").Value = '  +1.71%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0725'
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''7.68'
$ws.Range("E20").Value = '  +3.65%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("D23").Value = '''10.10'
$ws.Range("E23").Value = '  +7.50%  '
$ws.Range("E24").Value = '  -3.91%  '
$ws.Range("D25").Value = '''150.53'
$ws.Range("E25").Value = '  +2.37%  '
$ws.Range("D26").Value = '''6.97'
$ws.Range("E26").Value = '  -1.35%  '
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '''15.66'
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").Value = '''3.31'
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").Value = '1.443.35'
$ws.Range("E34").Value = '  +1.51%  '
$ws.Range("E35").Value = '  +1.07%  '
$ws.Range("D36").Value = '''2.34'
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("D37").Value = '''0.570'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").Value = '''0.884'
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").Value = '''0.871'
$ws.Range("E40").Value = '  +10.03%  '
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '''5.59'
$ws.Range("E43").Value = '  +0.81%  '
$ws.Range("D44").Value = '''66.68'
$ws.Range("E44").Value = '  +2.39%  '
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("D47").Value = '1.792.18'
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("D49").Value = '''86.35'
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").Value = '0.0₆0107'
$ws.Range("E50").Value = '  +2.41%  '
$ws.Range("D51").Value = '''0.0995'
$ws.Range("E51").Value = '  -1.53%  '
